$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H58").Value = 4186.6924
$ws.Range("J58").Value = 7428.4287
$ws.Range("L58").Value = 22285.2861
$ws.Range("N58").Value = -22585.2861

$ws.Range("H106").Value = 6796.8
$ws.Range("I106").Value = 4924.5
$ws.Range("K106").Value = 4924.5
$ws.Range("M106").Value = -4293.5

$ws.Range("H132").Value = 5582.603
$ws.Range("I132").Value = 5685.783
$ws.Range("K132").Value = 17057.349
$ws.Range("M132").Value = -14527.349

$ws.Range("H138").Value = 7263.1665
$ws.Range("I138").Value = 9212
$ws.Range("J138").Value = 5871.143
$ws.Range("K138").Value = 27636
$ws.Range("L138").Value = 17613.429
$ws.Range("M138").Value = -22496
$ws.Range("N138").Value = -27893.429

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 24391906
$ws.Range("I97").Value = 1409.4054
$ws.Range("K97").Value = 1409.4054
$ws.Range("M97").Value = -913.4054000000001

$ws.Range("H102").Value = 1627.6923
$ws.Range("I102").Value = 1596.6666
$ws.Range("K102").Value = 1596.6666
$ws.Range("M102").Value = 25.33339999999998

$ws.Range("H132").Value = 1318271.4
$ws.Range("I132").Value = 1564659.8
$ws.Range("J132").Value = 4200
$ws.Range("K132").Value = 4693979.4
$ws.Range("L132").Value = 12600
$ws.Range("M132").Value = -4691449.4
$ws.Range("N132").Value = -17660

$ws.Range("H133").Value = 72414
$ws.Range("J133").Value = 72414
$ws.Range("L133").Value = 72414
$ws.Range("N133").Value = -77474

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H13").Value = 74856
$ws.Range("J13").Value = 74856
$ws.Range("L13").Value = 74856
$ws.Range("N13").Value = -75192

$ws.Range("H86").Value = 4633.4287
$ws.Range("I86").Value = 3741.25
$ws.Range("K86").Value = 3741.25
$ws.Range("M86").Value = -2618.25

$ws.Range("H89").Value = 4633.4287
$ws.Range("I89").Value = 3741.25
$ws.Range("K89").Value = 18706.25
$ws.Range("M89").Value = -13090.25

$ws.Range("H94").Value = 1678.4
$ws.Range("I94").Value = 778.6
$ws.Range("K94").Value = 778.6
$ws.Range("M94").Value = -327.6

$ws.Range("H100").Value = 50000
$ws.Range("J100").Value = 50000
$ws.Range("L100").Value = 50000
$ws.Range("N100").Value = -52164

$ws.Range("H105").Value = 3751.3333
$ws.Range("I105").Value = 1493.4
$ws.Range("J105").Value = 6573.75
$ws.Range("K105").Value = 1493.4
$ws.Range("L105").Value = 6573.75
$ws.Range("M105").Value = 253.5999999999999
$ws.Range("N105").Value = -10067.75

$ws.Range("H109").Value = 64910.5
$ws.Range("J109").Value = 70000
$ws.Range("L109").Value = 70000
$ws.Range("N109").Value = -72774

$ws.Range("H132").Value = 96970.14
$ws.Range("J132").Value = 96997
$ws.Range("L132").Value = 96997
$ws.Range("N132").Value = -107117

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 17949.928
$ws.Range("I7").Value = 52675.684
$ws.Range("K7").Value = 52675.684
$ws.Range("M7").Value = -52562.684

$ws.Range("H107").Value = 644.1539
$ws.Range("I107").Value = 531.1667
$ws.Range("K107").Value = 531.1667
$ws.Range("M107").Value = 1388.8333

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 3837.6
$ws.Range("I7").Value = 4747
$ws.Range("J7").Value = 200
$ws.Range("K7").Value = 14241
$ws.Range("L7").Value = 600
$ws.Range("M7").Value = -14129
$ws.Range("N7").Value = -824

$ws.Range("H56").Value = 10000
$ws.Range("I56").Value = 10000
$ws.Range("K56").Value = 10000
$ws.Range("M56").Value = -9470

$ws.Range("H117").Value = 2291
$ws.Range("I117").Value = 802.7143
$ws.Range("K117").Value = 2408.1429
$ws.Range("M117").Value = 1033.8571

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3121.1538
$ws.Range("I22").Value = 2219.75
$ws.Range("J22").Value = 3521.7778
$ws.Range("K22").Value = 2219.75
$ws.Range("L22").Value = 3521.7778
$ws.Range("M22").Value = -1924.75
$ws.Range("N22").Value = -4111.7778

$ws.Range("H27").Value = 3121.1538
$ws.Range("I27").Value = 2219.75
$ws.Range("J27").Value = 3521.7778
$ws.Range("K27").Value = 2219.75
$ws.Range("L27").Value = 3521.7778
$ws.Range("M27").Value = -2112.75
$ws.Range("N27").Value = -3735.7778

$ws.Range("H61").Value = 15529.866
$ws.Range("I61").Value = 16564.143
$ws.Range("J61").Value = 1050
$ws.Range("K61").Value = 16564.143
$ws.Range("L61").Value = 1050
$ws.Range("M61").Value = -16362.143
$ws.Range("N61").Value = -1454

$ws.Range("H64").Value = 23465
$ws.Range("J64").Value = 23465
$ws.Range("L64").Value = 23465
$ws.Range("N64").Value = -23915

$ws.Range("H67").Value = 23465
$ws.Range("J67").Value = 23465
$ws.Range("L67").Value = 23465
$ws.Range("N67").Value = -25025

$ws.Range("H113").Value = 15529.866
$ws.Range("I113").Value = 16564.143
$ws.Range("J113").Value = 1050
$ws.Range("K113").Value = 16564.143
$ws.Range("L113").Value = 1050
$ws.Range("M113").Value = -14394.143
$ws.Range("N113").Value = -5390

$ws.Range("H133").Value = 74981.836
$ws.Range("J133").Value = 74981.836
$ws.Range("L133").Value = 74981.836
$ws.Range("N133").Value = -80041.836

$ws.Range("H136").Value = 13167661
$ws.Range("I136").Value = 9620624
$ws.Range("K136").Value = 28861872
$ws.Range("M136").Value = -28859322

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H63").Value = 0
$ws.Range("I63").Value = 0
$ws.Range("K63").Value = 0
$ws.Range("M63").ClearContents()

$ws.Range("H66").Value = 0
$ws.Range("I66").Value = 0
$ws.Range("K66").Value = 0
$ws.Range("M66").ClearContents()

$ws.Range("H81").Value = 3103.5625
$ws.Range("I81").Value = 2260.7144
$ws.Range("J81").Value = 9003.5
$ws.Range("K81").Value = 4521.4288
$ws.Range("L81").Value = 18007
$ws.Range("M81").Value = -3460.4288
$ws.Range("N81").Value = -20129

$ws.Range("H84").Value = 3103.5625
$ws.Range("I84").Value = 2260.7144
$ws.Range("J84").Value = 9003.5
$ws.Range("K84").Value = 22607.144
$ws.Range("L84").Value = 90035
$ws.Range("M84").Value = -17303.144
$ws.Range("N84").Value = -100643

$ws.Range("H107").Value = 30593.916
$ws.Range("I107").Value = 4839.8
$ws.Range("K107").Value = 14519.4
$ws.Range("M107").Value = -12599.4

$ws.Range("H141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("L141").Value = 0
